$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.018.48'
$ws.Range('E2').Value = '  +5.66%  '
$ws.Range('D3').Value = '2.261.45'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''302.13'
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').Value = '''92.86'
$ws.Range('E6').Value = '  +7.35%  '
$ws.Range('D7').Value = '''0.532'
$ws.Range('E7').Value = '  +3.90%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.485'
$ws.Range('E9').Value = '  +4.19%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '''54.88'
$ws.Range('E10').Value = '  +9.90%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '''32.77'
$ws.Range('E11').Value = '  +8.81%  '
$ws.Range('D12').Value = '''0.0799'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('E13').Value = '  +3.46%  '
$ws.Range('D14').Value = '''6.70'
$ws.Range('E14').Value = '  +4.62%  '
$ws.Range('D15').Value = '2.614.16'
$ws.Range('E15').Value = '  +2.38%  '
$ws.Range('D16').Value = '''14.14'
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('D17').Value = '2.279.24'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '''0.758'
$ws.Range('E18').Value = '  +4.26%  '
$ws.Range('D19').Value = '41.932.48'
$ws.Range('E19').Value = '  +5.48%  '
$ws.Range('D20').Value = '''12.16'
$ws.Range('E20').Value = '  +9.40%  '
$ws.Range('D21').Value = '0.0₃0907'
$ws.Range('E21').Value = '  +2.87%  '
$ws.Range('D22').Value = '''5.94'
$ws.Range('E22').Value = '  +3.93%  '
$ws.Range('D23').Value = '''67.28'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('D24').Value = '''241.89'
$ws.Range('E25').Value = '  +5.92%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''1.91'
$ws.Range('E27').Value = '  +4.76%  '
$ws.Range('D28').Value = '''23.91'
$ws.Range('E28').Value = '  +3.46%  '
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').Value = '''9.70'
$ws.Range('E30').Value = '  +5.72%  '
$ws.Range('D31').Value = '''34.23'
$ws.Range('E31').Value = '  +8.25%  '
$ws.Range('D32').Value = '''158.43'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '''5.17'
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('D35').Value = '''0.0742'
$ws.Range('E35').Value = '  +5.12%  '
$ws.Range('D36').Value = '''3.08'
$ws.Range('E36').Value = '  +6.03%  '
$ws.Range('D37').Value = '''2.41'
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('D38').Value = '''0.104'
$ws.Range('E38').Value = '  +6.76%  '
$ws.Range('D39').Value = '''16.57'
$ws.Range('E39').Value = '  +9.65%  '
$ws.Range('E40').Value = '  +4.58%  '
$ws.Range('D41').Value = '''1.80'
$ws.Range('E41').Value = '  +6.01%  '
$ws.Range('D42').Value = '''3.96'
$ws.Range('E42').Value = '  +7.03%  '
$ws.Range('D43').Value = '''20.18'
$ws.Range('E43').Value = '  +14.00%  '
$ws.Range('D44').Value = '2.054.88'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').Value = '''0.0280'
$ws.Range('E45').Value = '  +4.61%  '
$ws.Range('D46').Value = '''10.10'
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('D47').Value = '''2.91'
$ws.Range('E47').Value = '  +8.39%  '
$ws.Range('E48').Value = '  -4.65%  '
$ws.Range('D49').Value = '2.486.00'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('E50').Value = '  +3.39%  '
$ws.Range('E51').Value = '  +4.74%  '
